# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows with both Price (D) and Volume(1h) (E) changes
$ws.Range("D2").Value = '60.957.13'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '3.395.60'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("D5").Value = '''570.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").Value = '''141.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.43%  '
$ws.Range("D7").Value = '3.395.90'
$ws.Range("E7").Value = '  -0.75%  '
$ws.Range("D10").Value = '''7.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.43%  '
$ws.Range("D12").Value = '''0.394'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("D13").Value = '3.975.61'
$ws.Range("E13").Value = '  -0.91%  '
$ws.Range("D14").Value = '''28.34'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.06%  '
$ws.Range("D15").Value = '''0.124'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.67%  '
$ws.Range("D18").Value = '61.046.37'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = '''6.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").Value = '''13.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.65%  '
$ws.Range("D21").Value = '''9.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.12%  '
$ws.Range("D22").Value = '''385.58'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.51%  '
$ws.Range("D23").Value = '''0.559'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("D24").Value = '''74.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.66%  '
$ws.Range("D26").Value = '''0.0000117'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.00%  '
$ws.Range("D27").Value = '3.537.07'
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("D28").Value = '''0.178'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '''7.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("D31").Value = '''7.97'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.29%  '
$ws.Range("D32").Value = '''2.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("D33").Value = '''1.42'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.37%  '
$ws.Range("D35").Value = '''23.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("D36").Value = '''7.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").Value = '''167.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("D38").Value = '3.427.08'
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("D39").Value = '''4.99'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.18%  '
$ws.Range("D41").Value = '''28.58'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.97%  '
$ws.Range("D42").Value = '''0.0776'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("D45").Value = '''42.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D46").Value = '''4.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("D47").Value = '''1.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.49%  '
$ws.Range("D48").Value = '''1.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("D49").Value = '2.490.87'
$ws.Range("E49").Value = '  -3.40%  '
$ws.Range("D50").Value = '''23.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("D51").Value = '''6.82'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.41%  '

# Rows with only a Volume(1h) (E) change
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E40").Value = '  -4.42%  '

# Rows where two coins traded rank/position with each other
# (Coin, Link, Price and Volume all change together)
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.402.55'
$ws.Range("E16").Value = '  -0.94%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.0000171'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '''1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '''0.780'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.26%  '
